# Apply updated crypto price/volume figures (and restore the Polygon / wstETH row order)
# to the "cryptos" worksheet, matching the latest scrape snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '26.529.45'
$ws.Range("E2").Value = '  -7.71%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.685.27'
$ws.Range("E3").Value = '  -6.58%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '''1.005'
$ws.Range("E4").Value = '  +0.22%  '

# Row 5: BNB
$ws.Range("D5").Value = '''216.75'
$ws.Range("E5").Value = '  -6.51%  '

# Row 6: USDC
$ws.Range("E6").Value = '  +0.18%  '

# Row 7: XRP
$ws.Range("D7").Value = '''0.4972'
$ws.Range("E7").Value = '  -16.29%  '

# Row 8: Cardano
$ws.Range("D8").Value = '''0.2610'
$ws.Range("E8").Value = '  -6.05%  '

# Row 9: Solana
$ws.Range("D9").Value = '''21.64'
$ws.Range("E9").Value = '  -7.55%  '

# Row 10: Dogecoin
$ws.Range("D10").Value = '''0.06139'
$ws.Range("E10").Value = '  -10.35%  '

# Row 11: TRON
$ws.Range("D11").Value = '''0.07272'

# Row 12: WrappedEther
$ws.Range("D12").Value = '1.649.57'
$ws.Range("E12").Value = '  -8.72%  '

# Row 13: Polkadot
$ws.Range("D13").Value = '''4.421'
$ws.Range("E13").Value = '  -6.54%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '''0.5728'
$ws.Range("E14").Value = '  -8.76%  '

# Row 15: Polygon
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '1.913.20'
$ws.Range("E15").Value = '  -6.74%  '

# Row 16: ShibaInu
$ws.Range("D16").Value = '''0.000008294'
$ws.Range("E16").Value = '  -10.70%  '

# Row 17: Litecoin
$ws.Range("D17").Value = '''64.65'
$ws.Range("E17").Value = '  -14.21%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '26.557.17'
$ws.Range("E18").Value = '  -7.52%  '

# Row 19: Uniswap
$ws.Range("D19").Value = '''5.004'
$ws.Range("E19").Value = '  -8.52%  '

# Row 20: Dai
$ws.Range("E20").Value = '  +0.22%  '

# Row 21: Avalanche
$ws.Range("D21").Value = '''10.70'
$ws.Range("E21").Value = '  -6.64%  '

# Row 22: BitcoinCash
$ws.Range("D22").Value = '''182.60'
$ws.Range("E22").Value = '  -13.41%  '

# Row 23: Chainlink
$ws.Range("D23").Value = '''6.166'
$ws.Range("E23").Value = '  -10.19%  '

# Row 24: BinanceUSD
$ws.Range("D24").Value = '''1.007'
$ws.Range("E24").Value = '  +0.28%  '

# Row 25: Monero
$ws.Range("D25").Value = '''144.51'
$ws.Range("E25").Value = '  -6.43%  '

# Row 26: Cosmos
$ws.Range("D26").Value = '''7.568'
$ws.Range("E26").Value = '  -3.64%  '

# Row 27: Stellar
$ws.Range("D27").Value = '''0.1131'
$ws.Range("E27").Value = '  -11.33%  '

# Row 28: EthereumClassic
$ws.Range("D28").Value = '''15.28'
$ws.Range("E28").Value = '  -6.96%  '

# Row 29: Toncoin
$ws.Range("D29").Value = '''1.318'
$ws.Range("E29").Value = '  -8.58%  '

# Row 30: Hedera
$ws.Range("D30").Value = '''0.05586'
$ws.Range("E30").Value = '  -9.77%  '

# Row 31: PancakeSwap
$ws.Range("D31").Value = '''1.319'
$ws.Range("E31").Value = '  -7.13%  '

# Row 32: Filecoin
$ws.Range("D32").Value = '''3.480'
$ws.Range("E32").Value = '  -7.92%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").Value = '''3.463'
$ws.Range("E33").Value = '  -7.63%  '

# Row 34: LidoDAOToken
$ws.Range("E34").Value = '  -4.44%  '

# Row 35: ARBITRUM
$ws.Range("D35").Value = '''1.006'
$ws.Range("E35").Value = '  -5.04%  '

# Row 36: HuobiToken
$ws.Range("D36").Value = '''2.374'
$ws.Range("E36").Value = '  -5.02%  '

# Row 37: ImmutableX
$ws.Range("D37").Value = '''0.5871'
$ws.Range("E37").Value = '  -8.57%  '

# Row 38: MXToken
$ws.Range("D38").Value = '''2.638'
$ws.Range("E38").Value = '  -3.18%  '

# Row 39: VeChain
$ws.Range("D39").Value = '''0.01582'
$ws.Range("E39").Value = '  -7.48%  '

# Row 40: Maker
$ws.Range("D40").Value = '1.074.45'
$ws.Range("E40").Value = '  -6.10%  '

# Row 41: FraxShare
$ws.Range("D41").Value = '''5.910'
$ws.Range("E41").Value = '  -8.09%  '

# Row 42: TrustWalletToken
$ws.Range("D42").Value = '''0.8490'
$ws.Range("E42").Value = '  -2.04%  '

# Row 43: PaxDollar
$ws.Range("E43").Value = '  -0.18%  '

# Row 44: Quant
$ws.Range("D44").Value = '''98.31'
$ws.Range("E44").Value = '  -2.40%  '

# Row 45: RocketPoolETH
$ws.Range("D45").Value = '1.842.57'
$ws.Range("E45").Value = '  -6.22%  '

# Row 46: Aave
$ws.Range("D46").Value = '''56.24'
$ws.Range("E46").Value = '  -7.22%  '

# Row 47: BabyDogeCoin
$ws.Range("E47").Value = '  -6.70%  '

# Row 48: Frax
$ws.Range("D48").Value = '''1.004'
$ws.Range("E48").Value = '  -0.44%  '

# Row 49: EnergySwap
$ws.Range("D49").Value = '''8.056'
$ws.Range("E49").Value = '  -3.87%  '

# Row 50: Mantle
$ws.Range("D50").Value = '''0.4330'

# Row 51: Cronos
$ws.Range("D51").Value = '''0.05206'
$ws.Range("E51").Value = '  -4.87%  '
